$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

# Columns A and B hold text values that look like a date / a single
# CJK character respectively. Force them to be stored as text (not
# auto-converted to a date serial) by temporarily using a text number
# format, then reset the style back to "Normal" so no extra style is
# left behind on the new row (matching the rest of the sheet's rows).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/03"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "金"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = 201
